$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# LOP (L) values for rows that previously had blank LOP cells.
$ws.Range("L2").Value = 4
$ws.Range("L3").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("L5").Value = 1
$ws.Range("L6").Value = 1
$ws.Range("L7").Value = 1
$ws.Range("L8").Value = 1
$ws.Range("L9").Value = 1

# Rows 10-34: fill in LOP, Month and Year which were previously blank.
$ws.Range("L10").Value = 3
$ws.Range("M10").Value = "April"
$ws.Range("N10").Value = 2019

$ws.Range("L11").Value = 5
$ws.Range("M11").Value = "June"
$ws.Range("N11").Value = 2020

$ws.Range("L12").Value = 3
$ws.Range("M12").Value = "April"
$ws.Range("N12").Value = 2019

$ws.Range("L13").Value = 3
$ws.Range("M13").Value = "August"
$ws.Range("N13").Value = 2020

$ws.Range("L14").Value = 3
$ws.Range("M14").Value = "August"
$ws.Range("N14").Value = 2018

$ws.Range("L15").Value = 3
$ws.Range("M15").Value = "August"
$ws.Range("N15").Value = 2021

$ws.Range("L16").Value = 3
$ws.Range("M16").Value = "August"
$ws.Range("N16").Value = 2020

$ws.Range("L17").Value = 3
$ws.Range("M17").Value = "June"
$ws.Range("N17").Value = 2018

$ws.Range("L18").Value = 3
$ws.Range("M18").Value = "June"
$ws.Range("N18").Value = 2019

$ws.Range("L19").Value = 3
$ws.Range("M19").Value = "June"
$ws.Range("N19").Value = 2020

$ws.Range("L20").Value = 3
$ws.Range("M20").Value = "June"
$ws.Range("N20").Value = 2021

$ws.Range("L21").Value = 3
$ws.Range("M21").Value = "December"
$ws.Range("N21").Value = 2021

$ws.Range("L22").Value = 3
$ws.Range("M22").Value = "December"
$ws.Range("N22").Value = 2021

$ws.Range("L23").Value = 3
$ws.Range("M23").Value = "December"
$ws.Range("N23").Value = 2021

$ws.Range("L24").Value = 3
$ws.Range("M24").Value = "December"
$ws.Range("N24").Value = 2021

$ws.Range("L25").Value = 3
$ws.Range("M25").Value = "December"
$ws.Range("N25").Value = 2021

$ws.Range("L26").Value = 3
$ws.Range("M26").Value = "January"
$ws.Range("N26").Value = 2021

$ws.Range("L27").Value = 3
$ws.Range("M27").Value = "January"
$ws.Range("N27").Value = 2021

$ws.Range("L28").Value = 3
$ws.Range("M28").Value = "January"
$ws.Range("N28").Value = 2023

$ws.Range("L29").Value = 3
$ws.Range("M29").Value = "January"
$ws.Range("N29").Value = 2023

$ws.Range("L30").Value = 3
$ws.Range("M30").Value = "January"
$ws.Range("N30").Value = 2023

$ws.Range("L31").Value = 3
$ws.Range("M31").Value = "January"
$ws.Range("N31").Value = 2023

$ws.Range("L32").Value = 3
$ws.Range("M32").Value = "January"
$ws.Range("N32").Value = 2023

$ws.Range("L33").Value = 3
$ws.Range("M33").Value = "January"
$ws.Range("N33").Value = 2023

$ws.Range("L34").Value = 3
$ws.Range("M34").Value = "January"
$ws.Range("N34").Value = 2023

# Reflect the author's final cell selection.
$ws.Range("M14").Select()
